# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 165 of the data table,
# pushing the existing rows 165-256 down to 166-257 (dimension grows from
# A1:R256 to A1:R257).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 165, shifting rows 165..256 down to 166..257
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A165").Value2 = 4
$ws.Range("B165").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C165").Value2 = "Los Lagos"
$ws.Range("D165").Value2 = 44572
$ws.Range("E165").Value2 = 10
$ws.Range("F165").Value2 = 100112008
$ws.Range("G165").Value2 = "Coliflor"
$ws.Range("H165").Value2 = "Sin especificar"
$ws.Range("I165").Value2 = "Primera"
$ws.Range("J165").Value2 = 750
$ws.Range("K165").Value2 = 1400
$ws.Range("L165").Value2 = 1400
$ws.Range("M165").Value2 = 1400
$ws.Range("N165").Value2 = "$/unidad"
$ws.Range("O165").Value2 = "Región Metropolitana"
$ws.Range("P165").Value2 = 1400
$ws.Range("Q165").Value2 = 1
$ws.Range("R165").Value2 = "Hortaliza"
